# Applies the "Updated cryptos list" price/volume refresh to the single worksheet.
# Source data cells are plain text (t="inlineStr" in the original file) holding
# price/percentage strings such as "27.475.98" or "  -5.26%  ". Assigning these
# via .Value directly would let Excel auto-detect pure-decimal-looking strings
# (e.g. "1.001", "0.9999") as numbers and mangle them (trailing zeros lost,
# floating point noise introduced). Set-TextValue forces text by flipping the
# cell to the "@" (Text) number format before the write, then restores the
# original Style object so no visible formatting changes leak into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '27.475.98'
Set-TextValue $ws.Range('E2') '  -5.26%  '
Set-TextValue $ws.Range('D3') '1.838.58'
Set-TextValue $ws.Range('E3') '  -4.22%  '
Set-TextValue $ws.Range('D4') '1.001'
Set-TextValue $ws.Range('E4') '  -0.43%  '
Set-TextValue $ws.Range('D5') '313.33'
Set-TextValue $ws.Range('E5') '  -3.67%  '
Set-TextValue $ws.Range('D6') '0.9992'
Set-TextValue $ws.Range('E6') '  -0.45%  '
Set-TextValue $ws.Range('E7') '  -7.63%  '
Set-TextValue $ws.Range('D8') '0.3636'
Set-TextValue $ws.Range('E8') '  -4.49%  '
Set-TextValue $ws.Range('D9') '44.13'
Set-TextValue $ws.Range('E9') '  -3.39%  '
Set-TextValue $ws.Range('D10') '0.07221'
Set-TextValue $ws.Range('E10') '  -6.75%  '
Set-TextValue $ws.Range('D11') '0.9026'
Set-TextValue $ws.Range('E11') '  -7.55%  '
Set-TextValue $ws.Range('D12') '20.58'
Set-TextValue $ws.Range('E12') '  -8.68%  '
Set-TextValue $ws.Range('D13') '1.816.38'
Set-TextValue $ws.Range('E13') '  -5.40%  '
Set-TextValue $ws.Range('D14') '6.572'
Set-TextValue $ws.Range('E14') '  -5.27%  '
Set-TextValue $ws.Range('D15') '5.328'
Set-TextValue $ws.Range('E15') '  -6.58%  '
Set-TextValue $ws.Range('D16') '0.06791'
Set-TextValue $ws.Range('E16') '  -3.25%  '
Set-TextValue $ws.Range('D17') '0.9999'
Set-TextValue $ws.Range('E17') '  -0.65%  '
Set-TextValue $ws.Range('D18') '77.47'
Set-TextValue $ws.Range('E18') '  -8.42%  '
Set-TextValue $ws.Range('D19') '0.000009013'
Set-TextValue $ws.Range('E19') '  -4.92%  '
Set-TextValue $ws.Range('D20') '0.9987'
Set-TextValue $ws.Range('E20') '  -0.51%  '
Set-TextValue $ws.Range('D21') '15.34'
Set-TextValue $ws.Range('E21') '  -7.76%  '
Set-TextValue $ws.Range('D22') '27.487.29'
Set-TextValue $ws.Range('E22') '  -5.27%  '
Set-TextValue $ws.Range('D23') '4.935'
Set-TextValue $ws.Range('E23') '  -7.61%  '
Set-TextValue $ws.Range('D24') '10.53'
Set-TextValue $ws.Range('E24') '  -4.42%  '
Set-TextValue $ws.Range('D25') '2.010.40'
Set-TextValue $ws.Range('E25') '  -6.60%  '
Set-TextValue $ws.Range('D26') '2.025'
Set-TextValue $ws.Range('E26') '  -1.80%  '
Set-TextValue $ws.Range('E27') '  -2.94%  '
Set-TextValue $ws.Range('D28') '18.16'
Set-TextValue $ws.Range('E28') '  -4.56%  '
Set-TextValue $ws.Range('D29') '5.246'
Set-TextValue $ws.Range('E29') '  -6.22%  '
Set-TextValue $ws.Range('D30') '110.75'
Set-TextValue $ws.Range('E30') '  -5.88%  '
Set-TextValue $ws.Range('D31') '1.656'
Set-TextValue $ws.Range('E31') '  -9.27%  '
Set-TextValue $ws.Range('D32') '0.08849'
Set-TextValue $ws.Range('E32') '  -5.01%  '
Set-TextValue $ws.Range('D33') '0.7755'
Set-TextValue $ws.Range('E33') '  -9.30%  '
Set-TextValue $ws.Range('D34') '4.509'
Set-TextValue $ws.Range('E34') '  -11.18%  '
Set-TextValue $ws.Range('E35') '  -4.16%  '
Set-TextValue $ws.Range('B36') 'Frax'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D36') '0.9991'
Set-TextValue $ws.Range('E36') '  -0.47%  '
Set-TextValue $ws.Range('B37') 'ARBITRUM'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D37') '1.072'
Set-TextValue $ws.Range('E37') '  -13.35%  '
Set-TextValue $ws.Range('D38') '0.05336'
Set-TextValue $ws.Range('E38') '  -5.92%  '
Set-TextValue $ws.Range('D39') '1.083'
Set-TextValue $ws.Range('E39') '  -5.80%  '
Set-TextValue $ws.Range('D40') '0.01926'
Set-TextValue $ws.Range('E40') '  -5.55%  '
Set-TextValue $ws.Range('D41') '2.941'
Set-TextValue $ws.Range('E41') '  -5.18%  '
Set-TextValue $ws.Range('B42') 'FraxShare'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '6.824'
Set-TextValue $ws.Range('E42') '  -7.85%  '
Set-TextValue $ws.Range('B43') 'TheSandbox'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D43') '0.5054'
Set-TextValue $ws.Range('E43') '  -7.68%  '
Set-TextValue $ws.Range('E44') '  -6.79%  '
Set-TextValue $ws.Range('D45') '0.06610'
Set-TextValue $ws.Range('E45') '  -4.36%  '
Set-TextValue $ws.Range('D46') '8.222'
Set-TextValue $ws.Range('E46') '  -11.89%  '
Set-TextValue $ws.Range('D47') '0.4717'
Set-TextValue $ws.Range('E47') '  -8.66%  '
Set-TextValue $ws.Range('D48') '104.95'
Set-TextValue $ws.Range('E48') '  -4.65%  '
Set-TextValue $ws.Range('D49') '10.22'
Set-TextValue $ws.Range('E49') '  -8.59%  '
Set-TextValue $ws.Range('D50') '0.9989'
Set-TextValue $ws.Range('E50') '  -0.48%  '
Set-TextValue $ws.Range('D51') '1.627'
Set-TextValue $ws.Range('E51') '  -7.11%  '

Write-Output "Updated cryptos list: applied 104 cell changes across rows 2-51."
